# incluindo conexão de base de dados
#
# - rename sheet "pedido" -> "Descricao"
# - add a new, empty sheet "Loja" after it, and make it the active/selected sheet
# - selection on "Descricao" becomes a single cell (G1) instead of the whole column
# - selection on "Loja" is G6

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Descricao"
$ws1.Range("G1").Select() | Out-Null

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Loja"
$ws2.Range("G6").Select() | Out-Null
